{"js": "// Replace each division-problem cell's text with its new value.\n// The document is a single table of \"NNN\u00f7N=\" style problems; every\n// old value is unique, so a body.search()+insertText(\"Replace\") per\n// pair is unambiguous.\nconst replacements = [\n  [\"686\u00f79=\", \"312\u00f74=\"],\n  [\"576\u00f74=\", \"130\u00f78=\"],\n  [\"253\u00f77=\", \"374\u00f75=\"],\n  [\"497\u00f72=\", \"408\u00f72=\"],\n  [\"470\u00f78=\", \"543\u00f75=\"],\n  [\"826\u00f76=\", \"178\u00f77=\"],\n  [\"161\u00f73=\", \"411\u00f78=\"],\n  [\"131\u00f72=\", \"516\u00f76=\"],\n  [\"788\u00f79=\", \"914\u00f78=\"],\n  [\"346\u00f77=\", \"558\u00f77=\"],\n  [\"598\u00f73=\", \"909\u00f79=\"],\n  [\"818\u00f74=\", \"825\u00f74=\"],\n  [\"224\u00f78=\", \"214\u00f76=\"],\n  [\"962\u00f76=\", \"690\u00f72=\"],\n  [\"480\u00f75=\", \"669\u00f72=\"],\n  [\"239\u00f79=\", \"294\u00f75=\"],\n  [\"836\u00f79=\", \"442\u00f75=\"],\n  [\"731\u00f73=\", \"762\u00f74=\"],\n  [\"782\u00f73=\", \"334\u00f73=\"],\n  [\"446\u00f73=\", \"330\u00f79=\"],\n  [\"779\u00f73=\", \"709\u00f79=\"],\n  [\"968\u00f78=\", \"520\u00f75=\"],\n  [\"207\u00f78=\", \"142\u00f79=\"],\n  [\"680\u00f73=\", \"253\u00f73=\"],\n  [\"629\u00f75=\", \"551\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each division-problem cell's text with its new value.\n# The document is a single table of \"NNN\u00f7N=\" style problems; every\n# old value is unique, so Find/Replace (wdReplaceAll, scoped to an\n# exact match) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"686\u00f79=\", \"312\u00f74=\"),\n    @(\"576\u00f74=\", \"130\u00f78=\"),\n    @(\"253\u00f77=\", \"374\u00f75=\"),\n    @(\"497\u00f72=\", \"408\u00f72=\"),\n    @(\"470\u00f78=\", \"543\u00f75=\"),\n    @(\"826\u00f76=\", \"178\u00f77=\"),\n    @(\"161\u00f73=\", \"411\u00f78=\"),\n    @(\"131\u00f72=\", \"516\u00f76=\"),\n    @(\"788\u00f79=\", \"914\u00f78=\"),\n    @(\"346\u00f77=\", \"558\u00f77=\"),\n    @(\"598\u00f73=\", \"909\u00f79=\"),\n    @(\"818\u00f74=\", \"825\u00f74=\"),\n    @(\"224\u00f78=\", \"214\u00f76=\"),\n    @(\"962\u00f76=\", \"690\u00f72=\"),\n    @(\"480\u00f75=\", \"669\u00f72=\"),\n    @(\"239\u00f79=\", \"294\u00f75=\"),\n    @(\"836\u00f79=\", \"442\u00f75=\"),\n    @(\"731\u00f73=\", \"762\u00f74=\"),\n    @(\"782\u00f73=\", \"334\u00f73=\"),\n    @(\"446\u00f73=\", \"330\u00f79=\"),\n    @(\"779\u00f73=\", \"709\u00f79=\"),\n    @(\"968\u00f78=\", \"520\u00f75=\"),\n    @(\"207\u00f78=\", \"142\u00f79=\"),\n    @(\"680\u00f73=\", \"253\u00f73=\"),\n    @(\"629\u00f75=\", \"551\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text '$oldText' to replace with '$newText'\"\n    }\n}\n"}
